$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.443.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.16%  "

$ws.Range("D3").Value = "'3.067.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.98%  "

$ws.Range("E4").Value = "  -0.57%  "

$ws.Range("D5").Value = "'590.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.29%  "

$ws.Range("D6").Value = "'153.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.25%  "

$ws.Range("D8").Value = "'0.538"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.04%  "

$ws.Range("D9").Value = "'3.065.44"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.96%  "

$ws.Range("D10").Value = "'0.156"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.48%  "

$ws.Range("D11").Value = "'5.89"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.81%  "

$ws.Range("D12").Value = "'0.451"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.46%  "

$ws.Range("D13").Value = "'0.0000238"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.66%  "

$ws.Range("D14").Value = "'36.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.06%  "

$ws.Range("E15").Value = "  +1.17%  "

$ws.Range("D16").Value = "'3.570.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.13%  "

$ws.Range("D17").Value = "'7.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.13%  "

$ws.Range("D18").Value = "'63.330.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.91%  "

$ws.Range("D19").Value = "'3.059.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.37%  "

$ws.Range("D20").Value = "'484.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.71%  "

$ws.Range("D21").Value = "'14.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.58%  "

$ws.Range("D22").Value = "'0.708"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.57%  "

$ws.Range("D23").Value = "'7.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.18%  "

$ws.Range("D24").Value = "'2.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.80%  "

$ws.Range("D25").Value = "'82.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.02%  "

$ws.Range("D26").Value = "'12.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.71%  "

$ws.Range("D27").Value = "'10.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.32%  "

$ws.Range("E28").Value = "  +0.41%  "

$ws.Range("D29").Value = "'7.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.53%  "

$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "'2.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.57%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'2.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.78%  "

$ws.Range("D32").Value = "'0.998"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.73%  "

$ws.Range("D33").Value = "'27.33"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.51%  "

$ws.Range("D34").Value = "'0.112"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.84%  "

$ws.Range("E35").Value = "  +0.99%  "

$ws.Range("D36").Value = "'0.0₃0824"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.84%  "

$ws.Range("D37").Value = "'6.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.55%  "

$ws.Range("D38").Value = "'3.26"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.86%  "

$ws.Range("D39").Value = "'2.22"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.70%  "

$ws.Range("D40").Value = "'9.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.42%  "

$ws.Range("D41").Value = "'50.51"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.54%  "

$ws.Range("D42").Value = "'440.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.51%  "

$ws.Range("D43").Value = "'0.288"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.05%  "

$ws.Range("E44").Value = "  +2.64%  "

$ws.Range("D45").Value = "'0.0363"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.26%  "

$ws.Range("D46").Value = "'2.826.82"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.78%  "

$ws.Range("D47").Value = "'39.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.03%  "

$ws.Range("D48").Value = "'132.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.71%  "

$ws.Range("D49").Value = "'25.47"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.01%  "

$ws.Range("E50").Value = "  -0.02%  "

$ws.Range("D51").Value = "'2.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.31%  "
